$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted above the current row 204, pushing
# the existing rows 204-206 down to 205-207.
$ws.Rows.Item(204).Insert()

$ws.Cells.Item(204, 1).Value = 4
$ws.Cells.Item(204, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(204, 3).Value = "Los Lagos"
$ws.Cells.Item(204, 4).Value = 44595
$ws.Cells.Item(204, 4).NumberFormat = $ws.Cells.Item(205, 4).NumberFormat
$ws.Cells.Item(204, 5).Value = 10
$ws.Cells.Item(204, 6).Value = 100112044
$ws.Cells.Item(204, 7).Value = "Perejil"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 80
$ws.Cells.Item(204, 11).Value = 5000
$ws.Cells.Item(204, 12).Value = 5000
$ws.Cells.Item(204, 13).Value = 5000
$ws.Cells.Item(204, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(204, 15).Value = "Región Metropolitana"
$ws.Cells.Item(204, 16).Value = 1667
$ws.Cells.Item(204, 17).Value = 3
$ws.Cells.Item(204, 18).Value = "Hortaliza"
